# Auto-generated Excel COM-interop script to apply Belias_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(20, 8).Value = 7163.143
$ws.Cells.Item(20, 9).Value = 7163.143
$ws.Cells.Item(20, 11).Value = 7163.143
$ws.Cells.Item(20, 13).Value = -6933.143

$ws.Cells.Item(35, 8).Value = 7163.143
$ws.Cells.Item(35, 9).Value = 7163.143
$ws.Cells.Item(35, 11).Value = 7163.143
$ws.Cells.Item(35, 13).Value = -6784.143

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(26, 8).Value = 6793.3335
$ws.Cells.Item(26, 9).Value = 880
$ws.Cells.Item(26, 10).Value = 9750
$ws.Cells.Item(26, 11).Value = 880
$ws.Cells.Item(26, 12).Value = 9750
$ws.Cells.Item(26, 13).Value = -550
$ws.Cells.Item(26, 14).Value = -10410

$ws.Cells.Item(32, 8).Value = 3977.1
$ws.Cells.Item(32, 9).Value = 3327.6135
$ws.Cells.Item(32, 10).Value = 8740
$ws.Cells.Item(32, 11).Value = 3327.6135
$ws.Cells.Item(32, 12).Value = 8740
$ws.Cells.Item(32, 13).Value = -3040.6135
$ws.Cells.Item(32, 14).Value = -9314

$ws.Cells.Item(74, 8).Value = 91910.37
$ws.Cells.Item(74, 9).Value = 126001.75
$ws.Cells.Item(74, 10).Value = 1000
$ws.Cells.Item(74, 11).Value = 126001.75
$ws.Cells.Item(74, 12).Value = 1000
$ws.Cells.Item(74, 13).Value = -125127.75
$ws.Cells.Item(74, 14).Value = -2748

$ws.Cells.Item(77, 8).Value = 91910.37
$ws.Cells.Item(77, 9).Value = 126001.75
$ws.Cells.Item(77, 10).Value = 1000
$ws.Cells.Item(77, 11).Value = 630008.75
$ws.Cells.Item(77, 12).Value = 5000
$ws.Cells.Item(77, 13).Value = -625640.75
$ws.Cells.Item(77, 14).Value = -13736

$ws.Cells.Item(88, 8).Value = 1812.375
$ws.Cells.Item(88, 9).Value = 1900
$ws.Cells.Item(88, 10).Value = 1759.8
$ws.Cells.Item(88, 11).Value = 1900
$ws.Cells.Item(88, 12).Value = 1759.8
$ws.Cells.Item(88, 13).Value = -1494
$ws.Cells.Item(88, 14).Value = -2571.8

$ws.Cells.Item(91, 8).Value = 1812.375
$ws.Cells.Item(91, 9).Value = 1900
$ws.Cells.Item(91, 10).Value = 1759.8
$ws.Cells.Item(91, 11).Value = 1900
$ws.Cells.Item(91, 12).Value = 1759.8
$ws.Cells.Item(91, 13).Value = -496
$ws.Cells.Item(91, 14).Value = -4567.8

$ws.Cells.Item(92, 8).Value = 30000
$ws.Cells.Item(92, 10).Value = 30000
$ws.Cells.Item(92, 12).Value = 30000
$ws.Cells.Item(92, 14).Value = -34992

$ws.Cells.Item(122, 8).Value = 1557.2142
$ws.Cells.Item(122, 9).Value = 1466.1428
$ws.Cells.Item(122, 11).Value = 4398.428400000001
$ws.Cells.Item(122, 13).Value = -1948.428400000001

$ws.Cells.Item(135, 8).Value = 30780
$ws.Cells.Item(135, 10).Value = 30780
$ws.Cells.Item(135, 12).Value = 30780
$ws.Cells.Item(135, 14).Value = -40920

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1740.5476
$ws.Cells.Item(86, 9).Value = 1555.0344
$ws.Cells.Item(86, 10).Value = 2154.3845
$ws.Cells.Item(86, 11).Value = 1555.0344
$ws.Cells.Item(86, 12).Value = 2154.3845
$ws.Cells.Item(86, 13).Value = -432.0344
$ws.Cells.Item(86, 14).Value = -4400.3845

$ws.Cells.Item(89, 8).Value = 1740.5476
$ws.Cells.Item(89, 9).Value = 1555.0344
$ws.Cells.Item(89, 10).Value = 2154.3845
$ws.Cells.Item(89, 11).Value = 7775.172
$ws.Cells.Item(89, 12).Value = 10771.9225
$ws.Cells.Item(89, 13).Value = -2159.172
$ws.Cells.Item(89, 14).Value = -22003.9225

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 20475.2
$ws.Cells.Item(6, 9).Value = 23441.412
$ws.Cells.Item(6, 10).Value = 3666.6667
$ws.Cells.Item(6, 11).Value = 23441.412
$ws.Cells.Item(6, 12).Value = 3666.6667
$ws.Cells.Item(6, 13).Value = -23328.412
$ws.Cells.Item(6, 14).Value = -3892.6667

$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 13).Value = $null

$ws.Cells.Item(19, 8).Value = 515.7143
$ws.Cells.Item(19, 9).Value = 515.7143
$ws.Cells.Item(19, 11).Value = 515.7143
$ws.Cells.Item(19, 13).Value = -345.7143

$ws.Cells.Item(24, 8).Value = 515.7143
$ws.Cells.Item(24, 9).Value = 515.7143
$ws.Cells.Item(24, 11).Value = 515.7143
$ws.Cells.Item(24, 13).Value = -345.7143

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(9, 8).Value = 1543.75
$ws.Cells.Item(9, 9).Value = 900
$ws.Cells.Item(9, 10).Value = 1758.3334
$ws.Cells.Item(9, 11).Value = 2700
$ws.Cells.Item(9, 12).Value = 5275.0002
$ws.Cells.Item(9, 13).Value = -2476
$ws.Cells.Item(9, 14).Value = -5723.0002

$ws.Cells.Item(15, 8).Value = 538.25
$ws.Cells.Item(15, 9).Value = 35.333332
$ws.Cells.Item(15, 10).Value = 840
$ws.Cells.Item(15, 11).Value = 105.999996
$ws.Cells.Item(15, 12).Value = 2520
$ws.Cells.Item(15, 13).Value = 34.000004
$ws.Cells.Item(15, 14).Value = -2800

$ws.Cells.Item(16, 8).Value = 924.5
$ws.Cells.Item(16, 9).Value = 550
$ws.Cells.Item(16, 10).Value = 1299
$ws.Cells.Item(16, 11).Value = 1650
$ws.Cells.Item(16, 12).Value = 3897
$ws.Cells.Item(16, 13).Value = -1477
$ws.Cells.Item(16, 14).Value = -4243

$ws.Cells.Item(25, 8).Value = 1894.4445
$ws.Cells.Item(25, 9).Value = 650
$ws.Cells.Item(25, 11).Value = 1950
$ws.Cells.Item(25, 13).Value = -1781

$ws.Cells.Item(29, 8).Value = 400
$ws.Cells.Item(29, 9).Value = 200
$ws.Cells.Item(29, 10).Value = 500
$ws.Cells.Item(29, 11).Value = 600
$ws.Cells.Item(29, 12).Value = 1500
$ws.Cells.Item(29, 13).Value = -323
$ws.Cells.Item(29, 14).Value = -2054

$ws.Cells.Item(30, 8).Value = 1894.4445
$ws.Cells.Item(30, 9).Value = 650
$ws.Cells.Item(30, 11).Value = 1950
$ws.Cells.Item(30, 13).Value = -1848

$ws.Cells.Item(31, 8).Value = 1935.8334
$ws.Cells.Item(31, 9).Value = 1658.8889
$ws.Cells.Item(31, 10).Value = 2766.6667
$ws.Cells.Item(31, 11).Value = 4976.6667
$ws.Cells.Item(31, 12).Value = 8300.000100000001
$ws.Cells.Item(31, 13).Value = -4688.6667
$ws.Cells.Item(31, 14).Value = -8876.000100000001

$ws.Cells.Item(38, 8).Value = 156.09091
$ws.Cells.Item(38, 9).Value = 15
$ws.Cells.Item(38, 10).Value = 209
$ws.Cells.Item(38, 11).Value = 45
$ws.Cells.Item(38, 12).Value = 627
$ws.Cells.Item(38, 13).Value = 302
$ws.Cells.Item(38, 14).Value = -1321

$ws.Cells.Item(107, 8).Value = 291.625
$ws.Cells.Item(107, 10).Value = 298.16666
$ws.Cells.Item(107, 12).Value = 894.4999799999999
$ws.Cells.Item(107, 14).Value = -4734.49998

$ws.Cells.Item(134, 8).Value = 5884.952
$ws.Cells.Item(134, 9).Value = 4506.923
$ws.Cells.Item(134, 10).Value = 8124.25
$ws.Cells.Item(134, 11).Value = 13520.769
$ws.Cells.Item(134, 12).Value = 24372.75
$ws.Cells.Item(134, 13).Value = -8450.769
$ws.Cells.Item(134, 14).Value = -34512.75

$ws.Cells.Item(139, 8).Value = 3409.2354
$ws.Cells.Item(139, 9).Value = 1818.5294
$ws.Cells.Item(139, 10).Value = 4999.9414
$ws.Cells.Item(139, 11).Value = 5455.5882
$ws.Cells.Item(139, 12).Value = 14999.8242
$ws.Cells.Item(139, 13).Value = -315.5882000000001
$ws.Cells.Item(139, 14).Value = -25279.8242

$ws.Cells.Item(140, 8).Value = 2624.7354
$ws.Cells.Item(140, 9).Value = 974.9167
$ws.Cells.Item(140, 10).Value = 3524.6365
$ws.Cells.Item(140, 11).Value = 2924.7501
$ws.Cells.Item(140, 12).Value = 10573.9095
$ws.Cells.Item(140, 13).Value = 2255.2499
$ws.Cells.Item(140, 14).Value = -20933.9095

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2372.85
$ws.Cells.Item(80, 10).Value = 2432
$ws.Cells.Item(80, 12).Value = 2432
$ws.Cells.Item(80, 14).Value = -4428

$ws.Cells.Item(83, 8).Value = 2372.85
$ws.Cells.Item(83, 10).Value = 2432
$ws.Cells.Item(83, 12).Value = 12160
$ws.Cells.Item(83, 14).Value = -22144

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 979.1667
$ws.Cells.Item(82, 10).Value = 979.1667
$ws.Cells.Item(82, 12).Value = 979.1667
$ws.Cells.Item(82, 14).Value = -1701.1667

$ws.Cells.Item(85, 8).Value = 979.1667
$ws.Cells.Item(85, 10).Value = 979.1667
$ws.Cells.Item(85, 12).Value = 979.1667
$ws.Cells.Item(85, 14).Value = -3475.1667

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(45, 8).Value = 3990.4
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 10).Value = 3990.4
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 12).Value = 3990.4
$ws.Cells.Item(45, 13).Value = $null
$ws.Cells.Item(45, 14).Value = -4972.4

$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 14).Value = $null

$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 14).Value = $null
